# Updated symbol list on Tue Dec 13 15:58:58 UTC 2022 with GitHub Actions
#
# The "Price" column (D) holds numeric-looking values that are stored as
# TEXT in the workbook (inlineStr cells), not real numbers. A leading
# apostrophe forces Excel/COM to keep the assigned value as text instead
# of silently coercing it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 - Row 14: price refresh only ---
$ws.Range("D2").Value = "'274.67"
$ws.Range("D3").Value = "'22.91"
$ws.Range("D4").Value = "'6.314"
$ws.Range("D5").Value = "'0.06222"
$ws.Range("D6").Value = "'3.649"
$ws.Range("D7").Value = "'6.648"
$ws.Range("D8").Value = "'1.391"
$ws.Range("D9").Value = "'0.8333"
$ws.Range("D10").Value = "'0.01378"
$ws.Range("D11").Value = "'0.1598"
$ws.Range("D12").Value = "'0.08343"
$ws.Range("D13").Value = "'0.03520"
$ws.Range("D14").Value = "'0.03191"

# --- Row 15 - Row 26: coin list reshuffled (new entrant + re-ranked rows) ---
$ws.Range("B15").Value = "ProBitToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D15").Value = "'0.1261"
$ws.Range("E15").Value = "14ProBitTokenPROB"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'4.064"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitMartToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D17").Value = "'0.09298"
$ws.Range("E17").Value = "16BitMartTokenBMX"

$ws.Range("B18").Value = "BitForexToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D18").Value = "'0.001668"
$ws.Range("E18").Value = "17BitForexTokenBF"

$ws.Range("B19").Value = "CoinExToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D19").Value = "'0.04733"
$ws.Range("E19").Value = "18CoinExTokenCET"

$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D20").Value = "'0.006325"
$ws.Range("E20").Value = "19TigerCashTCH"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "'0.005702"
$ws.Range("E21").Value = "20HotbitTokenHTB"

$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").Value = "'0.001079"
$ws.Range("E22").Value = "21BitKanKAN"

$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").Value = "'0.0001502"
$ws.Range("E23").Value = "22NitroExNTX"

$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "'3.732"
$ws.Range("E24").Value = "23LEOLEO"

$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "'2.325"
$ws.Range("E25").Value = "24BTSETokenBTSE"

$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value = "'0.3353"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"

# --- scattered single-cell price refreshes ---
$ws.Range("D28").Value = "'0.0002707"
$ws.Range("D40").Value = "'0.04733"
$ws.Range("D41").Value = "'0.007118"

# --- Row 42 / Row 43: CEJI / BKEXToken swap ---
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003905"
$ws.Range("E42").Value = "41CEJICEJIWorstin24h"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1164"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# --- remaining scattered single-cell price refreshes ---
$ws.Range("D44").Value = "'0.01186"
$ws.Range("D45").Value = "'0.00006061"
$ws.Range("D46").Value = "'0.0009913"

$ws.Range("D48").Value = "'0.7831"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOIN"

$ws.Range("D49").Value = "'0.002398"
$ws.Range("D50").Value = "'0.00002403"
$ws.Range("D51").Value = "'0.01242"
